# Apply the changes described by the diff:
# 1. Update cell I2 (respFirstName) value from "Krishnaveni - Auto1" to "Krish - Auto1"
# 2. Move the active cell / selection on the worksheet from F9 to G11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdata")

# 1. Update the shared-string cell value
$ws.Range("I2").Value = "Krish - Auto1"

# 2. Change the current selection/active cell to G11 (bottom-right pane)
$ws.Range("G11").Select()
